$d = $word.ActiveDocument

# Update the header date
$d.Content.Find.Execute("2026-01-09 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-10 Saturday", 2)

# Helper to replace text inside a single table cell without disturbing the
# trailing paragraph mark / cell mark that Range.Text normally includes.
function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters from the end
    $r.End = $r.End - 1
    $r.Text = $newText
}

$table = $d.Tables.Item(1)

# Row 1 (grid row index 1)
Set-CellText $table 1 1 "57÷7="
Set-CellText $table 1 2 "20÷7="
Set-CellText $table 1 3 "34÷5="
Set-CellText $table 1 4 "91÷9="
Set-CellText $table 1 5 "14÷8="

# Row 5 (grid row index 5)
Set-CellText $table 5 1 "36÷2="
Set-CellText $table 5 2 "49÷5="
Set-CellText $table 5 3 "96÷4="
Set-CellText $table 5 4 "53÷8="
Set-CellText $table 5 5 "67÷6="

# Row 9 (grid row index 9)
Set-CellText $table 9 1 "16÷4="
Set-CellText $table 9 2 "63÷6="
Set-CellText $table 9 3 "35÷4="
Set-CellText $table 9 4 "83÷8="
Set-CellText $table 9 5 "82÷4="

# Row 13 (grid row index 13)
Set-CellText $table 13 1 "20÷4="
Set-CellText $table 13 2 "73÷8="
Set-CellText $table 13 3 "94÷5="
Set-CellText $table 13 4 "20÷6="
Set-CellText $table 13 5 "86÷4="

# Row 17 (grid row index 17)
Set-CellText $table 17 1 "73÷3="
Set-CellText $table 17 2 "44÷7="
Set-CellText $table 17 3 "72÷9="
Set-CellText $table 17 4 "50÷4="
Set-CellText $table 17 5 "97÷7="
